# Refresh market-price / leve-profit columns (H:N) for the rows that changed
# in this run of the scheduled Sheets updater, one worksheet
# (ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR) at a time.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets("ALC")
$ws.Range("H40").Value = 4104.1816
$ws.Range("I40").Value = 2583
$ws.Range("J40").Value = 4674.625
$ws.Range("K40").Value = 2583
$ws.Range("L40").Value = 4674.625
$ws.Range("M40").Value = -2408
$ws.Range("N40").Value = -5024.625

$ws.Range("H43").Value = 2775.25
$ws.Range("I43").Value = 2136.3333
$ws.Range("K43").Value = 2136.3333
$ws.Range("M43").Value = -2067.3333

$ws.Range("H64").Value = 5500
$ws.Range("J64").Value = 5500
$ws.Range("L64").Value = 5500
$ws.Range("N64").Value = -5996

$ws.Range("H67").Value = 5500
$ws.Range("J67").Value = 5500
$ws.Range("L67").Value = 5500
$ws.Range("N67").Value = -7216

$ws.Range("H99").Value = 3865
$ws.Range("I99").Value = 159
$ws.Range("J99").Value = 6644.5
$ws.Range("K99").Value = 477
$ws.Range("L99").Value = 19933.5
$ws.Range("M99").Value = 1021
$ws.Range("N99").Value = -22929.5

$ws.Range("H111").Value = 1110.1111
$ws.Range("I111").Value = 1123.875
$ws.Range("K111").Value = 3371.625
$ws.Range("M111").Value = -304.625

$ws.Range("H132").Value = 5676.489
$ws.Range("I132").Value = 4124.5127
$ws.Range("K132").Value = 12373.5381
$ws.Range("M132").Value = -9843.538100000002

$ws.Range("H137").Value = 1210.8667
$ws.Range("I137").Value = 828
$ws.Range("J137").Value = 3699.5
$ws.Range("K137").Value = 2484
$ws.Range("L137").Value = 11098.5
$ws.Range("M137").Value = 66
$ws.Range("N137").Value = -16198.5

$ws.Range("H139").Value = 107999.336
$ws.Range("J139").Value = 107999.336
$ws.Range("L139").Value = 107999.336
$ws.Range("N139").Value = -118279.336

$ws = $wb.Worksheets("ARM")
$ws.Range("H45").Value = 2669.3948
$ws.Range("J45").Value = 3828.2856
$ws.Range("L45").Value = 3828.2856
$ws.Range("N45").Value = -4582.2856

$ws.Range("H97").Value = 5822.143
$ws.Range("I97").Value = 5761.25
$ws.Range("K97").Value = 5761.25
$ws.Range("M97").Value = -5265.25

$ws.Range("H122").Value = 1797.5333
$ws.Range("I122").Value = 1497.3572
$ws.Range("K122").Value = 4492.071599999999
$ws.Range("M122").Value = -2042.071599999999

$ws = $wb.Worksheets("BSM")
$ws.Range("H22").Value = 708.8
$ws.Range("I22").Value = 708.8
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 708.8
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -535.8
$ws.Range("N22").Value = ""

$ws.Range("H99").Value = 30780.525
$ws.Range("I99").Value = 103720.3
$ws.Range("K99").Value = 103720.3
$ws.Range("M99").Value = -102222.3

$ws = $wb.Worksheets("CRP")
$ws.Range("H7").Value = 76.22727
$ws.Range("I7").Value = 58.307693
$ws.Range("K7").Value = 58.307693
$ws.Range("M7").Value = 54.692307

$ws.Range("H16").Value = 1079.8
$ws.Range("I16").Value = 1079.8
$ws.Range("K16").Value = 1079.8
$ws.Range("M16").Value = -792.8

$ws.Range("H17").Value = 6982.3335
$ws.Range("I17").Value = 7178.8
$ws.Range("K17").Value = 7178.8
$ws.Range("M17").Value = -7004.8

$ws.Range("H113").Value = 1079.8
$ws.Range("I113").Value = 1079.8
$ws.Range("K113").Value = 1079.8
$ws.Range("M113").Value = 1090.2

$ws = $wb.Worksheets("CUL")
$ws.Range("H7").Value = 5555
$ws.Range("I7").Value = 0
$ws.Range("K7").Value = 0
$ws.Range("M7").Value = ""

$ws.Range("H39").Value = 9420.714
$ws.Range("J39").Value = 9420.714
$ws.Range("L39").Value = 28262.142
$ws.Range("N39").Value = -28850.142

$ws.Range("H46").Value = 1258.091
$ws.Range("I46").Value = 604.375
$ws.Range("J46").Value = 3001.3333
$ws.Range("K46").Value = 1813.125
$ws.Range("L46").Value = 9003.999899999999
$ws.Range("M46").Value = -1722.125
$ws.Range("N46").Value = -9185.999899999999

$ws.Range("H128").Value = 121746
$ws.Range("I128").Value = 121746
$ws.Range("K128").Value = 365238
$ws.Range("M128").Value = -360258

$ws.Range("H131").Value = 5901774
$ws.Range("J131").Value = 9110566
$ws.Range("L131").Value = 27331698
$ws.Range("N131").Value = -27341778

$ws = $wb.Worksheets("GSM")
$ws.Range("H126").Value = 7854
$ws.Range("J126").Value = 8225
$ws.Range("L126").Value = 24675
$ws.Range("N126").Value = -29615

$ws = $wb.Worksheets("LTW")
$ws.Range("H7").Value = 4154.375
$ws.Range("I7").Value = 3979.7273
$ws.Range("J7").Value = 4538.6
$ws.Range("K7").Value = 3979.7273
$ws.Range("L7").Value = 4538.6
$ws.Range("M7").Value = -3867.7273
$ws.Range("N7").Value = -4762.6

$ws.Range("H40").Value = 2340.75
$ws.Range("I40").Value = 2298.7727
$ws.Range("J40").Value = 2802.5
$ws.Range("K40").Value = 2298.7727
$ws.Range("L40").Value = 2802.5
$ws.Range("M40").Value = -2162.7727
$ws.Range("N40").Value = -3074.5

$ws.Range("H46").Value = 7334.5
$ws.Range("I46").Value = 9516.909
$ws.Range("K46").Value = 9516.909
$ws.Range("M46").Value = -9328.909

$ws.Range("H93").Value = 1882.5358
$ws.Range("I93").Value = 1055.7778
$ws.Range("J93").Value = 2274.158
$ws.Range("K93").Value = 1055.7778
$ws.Range("L93").Value = 2274.158
$ws.Range("M93").Value = 192.2221999999999
$ws.Range("N93").Value = -4770.157999999999

$ws.Range("H100").Value = 4076.3845
$ws.Range("I100").Value = 3284.8572
$ws.Range("K100").Value = 3284.8572
$ws.Range("M100").Value = -2743.8572

$ws.Range("H126").Value = 4154.375
$ws.Range("I126").Value = 3979.7273
$ws.Range("J126").Value = 4538.6
$ws.Range("K126").Value = 11939.1819
$ws.Range("L126").Value = 13615.8
$ws.Range("M126").Value = -9469.1819
$ws.Range("N126").Value = -18555.8

$ws = $wb.Worksheets("WVR")
$ws.Range("H107").Value = 778.25
$ws.Range("I107").Value = 762.7143
$ws.Range("K107").Value = 2288.1429
$ws.Range("M107").Value = -368.1428999999998

$ws.Range("H122").Value = 2262.9666
$ws.Range("I122").Value = 2300.92
$ws.Range("K122").Value = 6902.76
$ws.Range("M122").Value = -4452.76

$ws.Range("H132").Value = 30273.324
$ws.Range("I132").Value = 31603.5
$ws.Range("J132").Value = 5000
$ws.Range("K132").Value = 94810.5
$ws.Range("L132").Value = 15000
$ws.Range("M132").Value = -92280.5
$ws.Range("N132").Value = -20060
